# Update "想去人数" (F column) values on sheets "展览" and "全部类型"
# to reflect refreshed data pulled at commit 456a3b4.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 6819
$wsExhibit.Range("F3").Value = 89
$wsExhibit.Range("F5").Value = 440
$wsExhibit.Range("F6").Value = 149
$wsExhibit.Range("F8").Value = 58
$wsExhibit.Range("F9").Value = 199
$wsExhibit.Range("F10").Value = 1280
$wsExhibit.Range("F12").Value = 106
$wsExhibit.Range("F14").Value = 132
$wsExhibit.Range("F15").Value = 18
$wsExhibit.Range("F16").Value = 379
$wsExhibit.Range("F17").Value = 45
$wsExhibit.Range("F19").Value = 4830
$wsExhibit.Range("F20").Value = 87
$wsExhibit.Range("F21").Value = 74
$wsExhibit.Range("F22").Value = 291
$wsExhibit.Range("F23").Value = 204

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F5").Value = 440
$wsAll.Range("F6").Value = 149
$wsAll.Range("F7").Value = 6480
$wsAll.Range("F8").Value = 58
$wsAll.Range("F9").Value = 199
$wsAll.Range("F10").Value = 1280
$wsAll.Range("F12").Value = 106
$wsAll.Range("F13").Value = 0
$wsAll.Range("F14").Value = 0
$wsAll.Range("F22").Value = 74
$wsAll.Range("F23").Value = 291
$wsAll.Range("F24").Value = 204
$wsAll.Range("F25").Value = 143
